$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.20568568574637
$ws.Range("C2").Value = 16.82946031282886
$ws.Range("D2").Value = 17.54084446684263

$ws.Range("B3").Value = 1.527678311055065
$ws.Range("C3").Value = 1.807776088457613
$ws.Range("D3").Value = 2.313279794398638

$ws.Range("B4").Value = 0.3240385233436823
$ws.Range("C4").Value = 0.3858151876220725
$ws.Range("D4").Value = 0.4983013320012676

$ws.Range("B5").Value = 79.76350511383396
$ws.Range("C5").Value = 80.6550332882285
$ws.Range("D5").Value = 81.64584448885297
